# Auto-generated Excel COM-interop script to apply scheduled market-data update
# to the FFXIV leve-profit tracker workbook (Sophia_Profits.xlsx).
# For each affected row, columns H-N (current market prices / computed profits)
# are refreshed with newly fetched values; a few cells are added or cleared
# where the corresponding metric is no longer (or newly) applicable.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value2 = 5500
$ws.Range("I40").Value2 = 5500
$ws.Range("K40").Value2 = 5500
$ws.Range("M40").Value2 = -5325

$ws.Range("H137").Value2 = 1829.8108
$ws.Range("I137").Value2 = 2016.0454
$ws.Range("J137").Value2 = 1556.6666
$ws.Range("K137").Value2 = 6048.1362
$ws.Range("L137").Value2 = 4669.9998
$ws.Range("M137").Value2 = -3498.1362
$ws.Range("N137").Value2 = -9769.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 2600.4
$ws.Range("I45").Value2 = 2600.4
$ws.Range("J45").Value2 = 0
$ws.Range("K45").Value2 = 2600.4
$ws.Range("L45").Value2 = 0
$ws.Range("M45").Value2 = -2223.4
$ws.Range("N45").ClearContents()

$ws.Range("H74").Value2 = 10821
$ws.Range("I74").Value2 = 11646.647
$ws.Range("J74").Value2 = 7312
$ws.Range("K74").Value2 = 11646.647
$ws.Range("L74").Value2 = 7312
$ws.Range("M74").Value2 = -10772.647
$ws.Range("N74").Value2 = -9060

$ws.Range("H77").Value2 = 10821
$ws.Range("I77").Value2 = 11646.647
$ws.Range("J77").Value2 = 7312
$ws.Range("K77").Value2 = 58233.235
$ws.Range("L77").Value2 = 36560
$ws.Range("M77").Value2 = -53865.235
$ws.Range("N77").Value2 = -45296

$ws.Range("H97").Value2 = 2500
$ws.Range("I97").Value2 = 2500
$ws.Range("K97").Value2 = 2500
$ws.Range("M97").Value2 = -2004

$ws.Range("H110").Value2 = 3817.182
$ws.Range("I110").Value2 = 1765.5555
$ws.Range("J110").Value2 = 13049.5
$ws.Range("K110").Value2 = 1765.5555
$ws.Range("L110").Value2 = 13049.5
$ws.Range("M110").Value2 = 279.4445000000001
$ws.Range("N110").Value2 = -17139.5

$ws.Range("H122").Value2 = 2874.25
$ws.Range("J122").Value2 = 0
$ws.Range("L122").Value2 = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 4061.8462
$ws.Range("I20").Value2 = 1164.1818
$ws.Range("J20").Value2 = 19999
$ws.Range("K20").Value2 = 1164.1818
$ws.Range("L20").Value2 = 19999
$ws.Range("M20").Value2 = -917.1818000000001
$ws.Range("N20").Value2 = -20493

$ws.Range("H51").Value2 = 0
$ws.Range("J51").Value2 = 0
$ws.Range("L51").Value2 = 0
$ws.Range("N51").ClearContents()

$ws.Range("H86").Value2 = 8099.75
$ws.Range("I86").Value2 = 2799.5
$ws.Range("J86").Value2 = 13400
$ws.Range("K86").Value2 = 2799.5
$ws.Range("L86").Value2 = 13400
$ws.Range("M86").Value2 = -1676.5
$ws.Range("N86").Value2 = -15646

$ws.Range("H89").Value2 = 8099.75
$ws.Range("I89").Value2 = 2799.5
$ws.Range("J89").Value2 = 13400
$ws.Range("K89").Value2 = 13997.5
$ws.Range("L89").Value2 = 67000
$ws.Range("M89").Value2 = -8381.5
$ws.Range("N89").Value2 = -78232

$ws.Range("H94").Value2 = 2189.5557
$ws.Range("I94").Value2 = 3301.4
$ws.Range("J94").Value2 = 799.75
$ws.Range("K94").Value2 = 3301.4
$ws.Range("L94").Value2 = 799.75
$ws.Range("M94").Value2 = -2850.4
$ws.Range("N94").Value2 = -1701.75

$ws.Range("H107").Value2 = 603.3333
$ws.Range("I107").Value2 = 595.8333
$ws.Range("J107").Value2 = 633.3333
$ws.Range("K107").Value2 = 595.8333
$ws.Range("L107").Value2 = 633.3333
$ws.Range("M107").Value2 = 1324.1667
$ws.Range("N107").Value2 = -4473.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1743.35
$ws.Range("I31").Value2 = 1573.5834
$ws.Range("J31").Value2 = 1998
$ws.Range("K31").Value2 = 1573.5834
$ws.Range("L31").Value2 = 1998
$ws.Range("M31").Value2 = -1278.5834
$ws.Range("N31").Value2 = -2588

$ws.Range("H34").Value2 = 1743.35
$ws.Range("I34").Value2 = 1573.5834
$ws.Range("J34").Value2 = 1998
$ws.Range("K34").Value2 = 1573.5834
$ws.Range("L34").Value2 = 1998
$ws.Range("M34").Value2 = -1371.5834
$ws.Range("N34").Value2 = -2402

$ws.Range("H99").Value2 = 2053.5
$ws.Range("I99").Value2 = 2053.5
$ws.Range("K99").Value2 = 2053.5
$ws.Range("M99").Value2 = -555.5

$ws.Range("H107").Value2 = 971.4
$ws.Range("I107").Value2 = 990.7857
$ws.Range("J107").Value2 = 700
$ws.Range("K107").Value2 = 990.7857
$ws.Range("L107").Value2 = 700
$ws.Range("M107").Value2 = 929.2143
$ws.Range("N107").Value2 = -4540

$ws.Range("H126").Value2 = 2053.5
$ws.Range("I126").Value2 = 2053.5
$ws.Range("K126").Value2 = 6160.5
$ws.Range("M126").Value2 = -3690.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value2 = 1500
$ws.Range("J25").Value2 = 1500
$ws.Range("L25").Value2 = 4500
$ws.Range("N25").Value2 = -4838

$ws.Range("H30").Value2 = 1500
$ws.Range("J30").Value2 = 1500
$ws.Range("L30").Value2 = 4500
$ws.Range("N30").Value2 = -4704

$ws.Range("H129").Value2 = 2068.3
$ws.Range("J129").Value2 = 2564
$ws.Range("L129").Value2 = 7692
$ws.Range("N129").Value2 = -17692

$ws.Range("H137").Value2 = 12979.8
$ws.Range("J137").Value2 = 24139.6
$ws.Range("L137").Value2 = 72418.79999999999
$ws.Range("N137").Value2 = -82618.79999999999

$ws.Range("H140").Value2 = 1922.3
$ws.Range("I140").Value2 = 1802.5555
$ws.Range("K140").Value2 = 5407.666499999999
$ws.Range("M140").Value2 = -227.6664999999994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 1648.75
$ws.Range("I97").Value2 = 1933.3334
$ws.Range("J97").Value2 = 795
$ws.Range("K97").Value2 = 1933.3334
$ws.Range("L97").Value2 = 795
$ws.Range("M97").Value2 = -1437.3334
$ws.Range("N97").Value2 = -1787

$ws.Range("H102").Value2 = 3745.75
$ws.Range("I102").Value2 = 3745.75
$ws.Range("K102").Value2 = 3745.75
$ws.Range("M102").Value2 = -2123.75

$ws.Range("H113").Value2 = 4833.3335
$ws.Range("I113").Value2 = 4500
$ws.Range("J113").Value2 = 5500
$ws.Range("K113").Value2 = 4500
$ws.Range("L113").Value2 = 5500
$ws.Range("M113").Value2 = -2330
$ws.Range("N113").Value2 = -9840

$ws.Range("H126").Value2 = 3037.25
$ws.Range("I126").Value2 = 3133.3333
$ws.Range("K126").Value2 = 9399.999899999999
$ws.Range("M126").Value2 = -6929.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 11864.071
$ws.Range("I22").Value2 = 15028.143
$ws.Range("J22").Value2 = 8700
$ws.Range("K22").Value2 = 15028.143
$ws.Range("L22").Value2 = 8700
$ws.Range("M22").Value2 = -14733.143
$ws.Range("N22").Value2 = -9290

$ws.Range("H27").Value2 = 11864.071
$ws.Range("I27").Value2 = 15028.143
$ws.Range("J27").Value2 = 8700
$ws.Range("K27").Value2 = 15028.143
$ws.Range("L27").Value2 = 8700
$ws.Range("M27").Value2 = -14921.143
$ws.Range("N27").Value2 = -8914

$ws.Range("H46").Value2 = 3083.3333
$ws.Range("I46").Value2 = 2700
$ws.Range("K46").Value2 = 2700
$ws.Range("M46").Value2 = -2512

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value2 = 891.3333
$ws.Range("I113").Value2 = 891.3333
$ws.Range("K113").Value2 = 2673.9999
$ws.Range("M113").Value2 = -503.9998999999998

$ws.Range("H126").Value2 = 1737.5
$ws.Range("I126").Value2 = 1850
$ws.Range("K126").Value2 = 5550
$ws.Range("M126").Value2 = -3080
